# Auto-generated edit script applying the Sargatanas_Profits.xlsx diff
# to the corresponding worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I2").Value = 38.94737
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 38.94737
$ws.Range("L2").Value = 60
$ws.Range("M2").Value = 74.05262999999999
$ws.Range("N2").Value = -286
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("H41").Value = 10419030
$ws.Range("I41").Value = 17858854
$ws.Range("J41").Value = 3278.8
$ws.Range("K41").Value = 17858854
$ws.Range("L41").Value = 3278.8
$ws.Range("M41").Value = -17858414
$ws.Range("N41").Value = -4158.8
$ws.Range("H76").Value = 7273.75
$ws.Range("I76").Value = 5741.857
$ws.Range("K76").Value = 5741.857
$ws.Range("M76").Value = -5426.857
$ws.Range("H79").Value = 7273.75
$ws.Range("I79").Value = 5741.857
$ws.Range("K79").Value = 5741.857
$ws.Range("M79").Value = -4649.857
$ws.Range("H116").Value = 20841416
$ws.Range("I116").Value = 62502500
$ws.Range("J116").Value = 10873.75
$ws.Range("K116").Value = 62502500
$ws.Range("L116").Value = 10873.75
$ws.Range("M116").Value = -62499058
$ws.Range("N116").Value = -17757.75
$ws.Range("H125").Value = 200001460
$ws.Range("I125").Value = 500000000
$ws.Range("K125").Value = 4500000000
$ws.Range("M125").Value = -4499997540
$ws.Range("H132").Value = 1151.2709
$ws.Range("I132").Value = 827.3721
$ws.Range("J132").Value = 3936.8
$ws.Range("K132").Value = 2482.1163
$ws.Range("L132").Value = 11810.4
$ws.Range("M132").Value = 47.88369999999986
$ws.Range("N132").Value = -16870.4
$ws.Range("H137").Value = 4568.0244
$ws.Range("J137").Value = 7592.8823
$ws.Range("L137").Value = 22778.6469
$ws.Range("N137").Value = -27878.6469
$ws.Range("H138").Value = 1728704.4
$ws.Range("I138").Value = 2533
$ws.Range("J138").Value = 2637215.8
$ws.Range("K138").Value = 7599
$ws.Range("L138").Value = 7911647.399999999
$ws.Range("M138").Value = -2459
$ws.Range("N138").Value = -7921927.399999999
$ws.Range("L21").ClearContents()
$ws.Range("L23").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1951.24
$ws.Range("I2").Value = 1800.0588
$ws.Range("J2").Value = 2272.5
$ws.Range("K2").Value = 1800.0588
$ws.Range("L2").Value = 2272.5
$ws.Range("M2").Value = -1687.0588
$ws.Range("N2").Value = -2498.5
$ws.Range("H61").Value = 37040024
$ws.Range("I61").Value = 2280.6667
$ws.Range("J61").Value = 166672110
$ws.Range("K61").Value = 2280.6667
$ws.Range("L61").Value = 166672110
$ws.Range("M61").Value = -2068.6667
$ws.Range("N61").Value = -166672534
$ws.Range("H102").Value = 4370.6
$ws.Range("I102").Value = 3710.182
$ws.Range("K102").Value = 3710.182
$ws.Range("M102").Value = -2088.182
$ws.Range("H116").Value = 1951.24
$ws.Range("I116").Value = 1800.0588
$ws.Range("J116").Value = 2272.5
$ws.Range("K116").Value = 1800.0588
$ws.Range("L116").Value = 2272.5
$ws.Range("M116").Value = 493.9412
$ws.Range("N116").Value = -6860.5
$ws.Range("H136").Value = 37040024
$ws.Range("I136").Value = 2280.6667
$ws.Range("J136").Value = 166672110
$ws.Range("K136").Value = 6842.000100000001
$ws.Range("L136").Value = 500016330
$ws.Range("M136").Value = -4292.000100000001
$ws.Range("N136").Value = -500021430

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1951.24
$ws.Range("I3").Value = 1800.0588
$ws.Range("J3").Value = 2272.5
$ws.Range("K3").Value = 1800.0588
$ws.Range("L3").Value = 2272.5
$ws.Range("M3").Value = -1686.0588
$ws.Range("N3").Value = -2500.5
$ws.Range("H22").Value = 275.4
$ws.Range("I22").Value = 269.5
$ws.Range("J22").Value = 299
$ws.Range("K22").Value = 269.5
$ws.Range("L22").Value = 299
$ws.Range("M22").Value = -96.5
$ws.Range("N22").Value = -645
$ws.Range("H99").Value = 18184018
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 30305530
$ws.Range("K99").Value = 1750
$ws.Range("L99").Value = 30305530
$ws.Range("M99").Value = -252
$ws.Range("N99").Value = -30308526
$ws.Range("H105").Value = 2771.8223
$ws.Range("I105").Value = 1866.4375
$ws.Range("K105").Value = 1866.4375
$ws.Range("M105").Value = -119.4375
$ws.Range("H107").Value = 125012890
$ws.Range("I107").Value = 187517010
$ws.Range("K107").Value = 187517010
$ws.Range("M107").Value = -187515090
$ws.Range("H134").Value = 6101445
$ws.Range("I134").Value = 8335750
$ws.Range("J134").Value = 7885.909
$ws.Range("K134").Value = 25007250
$ws.Range("L134").Value = 23657.727
$ws.Range("M134").Value = -25004715
$ws.Range("N134").Value = -28727.727
$ws.Range("H135").Value = 99780
$ws.Range("J135").Value = 99780
$ws.Range("L135").Value = 99780
$ws.Range("N135").Value = -109920

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4896.1787
$ws.Range("I132").Value = 3343.5
$ws.Range("K132").Value = 10030.5
$ws.Range("M132").Value = -7500.5
$ws.Range("H134").Value = 5321.5
$ws.Range("I134").Value = 4200
$ws.Range("J134").Value = 5616.6313
$ws.Range("K134").Value = 12600
$ws.Range("L134").Value = 16849.8939
$ws.Range("M134").Value = -10065
$ws.Range("N134").Value = -21919.8939
$ws.Range("H137").Value = 74999.5
$ws.Range("J137").Value = 75000
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("L37").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1193.0588
$ws.Range("I97").Value = 1025.4667
$ws.Range("J97").Value = 2450
$ws.Range("K97").Value = 1025.4667
$ws.Range("L97").Value = 2450
$ws.Range("M97").Value = -529.4666999999999
$ws.Range("N97").Value = -3442
$ws.Range("H107").Value = 421236.06
$ws.Range("I107").Value = 800173.8
$ws.Range("J107").Value = 194.11111
$ws.Range("K107").Value = 800173.8
$ws.Range("L107").Value = 194.11111
$ws.Range("M107").Value = -798253.8
$ws.Range("N107").Value = -4034.11111
$ws.Range("H132").Value = 3105
$ws.Range("I132").Value = 2141.8333
$ws.Range("K132").Value = 6425.499899999999
$ws.Range("M132").Value = -3895.499899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1346.1428
$ws.Range("J22").Value = 2950.4
$ws.Range("L22").Value = 2950.4
$ws.Range("N22").Value = -3540.4
$ws.Range("H27").Value = 1346.1428
$ws.Range("J27").Value = 2950.4
$ws.Range("L27").Value = 2950.4
$ws.Range("N27").Value = -3164.4
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("H55").Value = 327.47827
$ws.Range("J55").Value = 467.30768
$ws.Range("L55").Value = 467.30768
$ws.Range("N55").Value = -813.30768
$ws.Range("H80").Value = 60000
$ws.Range("J80").Value = 60000
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -62246
$ws.Range("H83").Value = 60000
$ws.Range("J83").Value = 60000
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -191232
$ws.Range("L36").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 17576308
$ws.Range("I81").Value = 1522884.9
$ws.Range("J81").Value = 40051100
$ws.Range("K81").Value = 3045769.8
$ws.Range("L81").Value = 80102200
$ws.Range("M81").Value = -3044708.8
$ws.Range("N81").Value = -80104322
$ws.Range("H84").Value = 17576308
$ws.Range("I84").Value = 1522884.9
$ws.Range("J84").Value = 40051100
$ws.Range("K84").Value = 15228849
$ws.Range("L84").Value = 400511000
$ws.Range("M84").Value = -15223545
$ws.Range("N84").Value = -400521608
$ws.Range("H107").Value = 37039016
$ws.Range("J107").Value = 55557190
$ws.Range("L107").Value = 166671570
$ws.Range("N107").Value = -166675410
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("H132").Value = 8898.4
$ws.Range("I132").Value = 9331.333000000001
$ws.Range("J132").Value = 8249
$ws.Range("K132").Value = 27993.999
$ws.Range("L132").Value = 24747
$ws.Range("M132").Value = -25463.999
$ws.Range("N132").Value = -29807
$ws.Range("H135").Value = 61500
$ws.Range("J135").Value = 61500
$ws.Range("L135").Value = 61500
$ws.Range("N135").Value = -71640
$ws.Range("H137").Value = 71249.75
$ws.Range("J137").Value = 71249.75
$ws.Range("L137").Value = 71249.75
$ws.Range("N137").Value = -81449.75
$ws.Range("M126").ClearContents()
